$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.136.93"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.92%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.264.86"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.10%  "

$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "308.06"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.27%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.45"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.67%  "

$ws.Range("E7").Value = "  -1.30%  "

$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.489"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.47%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.30"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.84%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0787"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.24%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.114"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.37%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.81"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.16%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.612.52"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.09%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.58"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.71%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.292.66"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.01%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.789"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.09%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "41.954.59"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.04%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.22"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.96%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0901"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.06%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.95"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.15%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.50"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.52%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.32"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.83%  "

$ws.Range("E24").Value = "  -1.35%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.97"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.49%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.02%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "23.54"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.02%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "36.53"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.29%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.13"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.17%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.52"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.98%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "164.08"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.88%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.23"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.02%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.01%  "

$ws.Range("E34").Value = "  +0.21%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0735"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.31%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.35"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.63%  "

$ws.Range("E37").Value = "  -0.02%  "

$ws.Range("E38").Value = "  -3.78%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.115"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.77%  "

$ws.Range("E40").Value = "  -3.69%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.14"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.03%  "

$ws.Range("E42").Value = "  -5.90%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.954.57"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.50%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "18.86"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.78%  "

$ws.Range("E45").Value = "  -2.04%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.92"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.42%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.78"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.70%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "53.36"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.35%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.484.93"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.98%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "92.23"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.15%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "71.43"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.81%  "
